# Updates the "cryptos" worksheet with refreshed Price (D) and Volume(1h) (E)
# figures, matching a GitHub Actions scheduled data refresh.
#
# Note: several Price values are plain decimals (e.g. "0.5063", "0.9930")
# which Excel's automatic type inference would otherwise coerce into numbers
# (dropping significant trailing zeros / exact text formatting). To force
# those specific values to be stored as literal text -- just like the
# original inline-string cells -- they are written with a leading
# apostrophe, Excel's standard "treat as text" quote-prefix marker. Values
# that are not parsed as pure numbers (e.g. "25.820.23", URLs, names,
# percentages with padding spaces) are written as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.820.23"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "1.629.14"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'215.48"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "'0.5063"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("D8").Value = "'0.2578"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "'0.06433"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").Value = "'19.36"
$ws.Range("E10").Value = "  -2.12%  "
$ws.Range("D11").Value = "'0.07802"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "'4.254"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "1.627.61"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").Value = "1.852.88"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "'0.5578"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").Value = "'63.25"
$ws.Range("E16").Value = "  -1.57%  "
$ws.Range("D17").Value = "0.0₅7518"
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").Value = "25.816.66"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "'193.05"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("D21").Value = "'4.297"
$ws.Range("E21").Value = "  -2.93%  "
$ws.Range("D22").Value = "'9.790"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").Value = "'1.004"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "'1.823"
$ws.Range("E25").Value = "  -3.41%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'140.15"
$ws.Range("E26").Value = "  -2.21%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.1264"
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("D28").Value = "'6.713"
$ws.Range("E28").Value = "  -2.13%  "
$ws.Range("D29").Value = "'15.40"
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("D30").Value = "'1.239"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").Value = "'0.04856"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").Value = "'3.275"
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "'3.179"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("D35").Value = "'2.379"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").Value = "'0.8935"
$ws.Range("E36").Value = "  -2.07%  "
$ws.Range("D37").Value = "1.135.04"
$ws.Range("E37").Value = "  +3.87%  "
$ws.Range("D38").Value = "'2.556"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("D39").Value = "'0.5463"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("D41").Value = "'0.9930"
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("D42").Value = "'5.560"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").Value = "'0.7934"
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("D44").Value = "'97.17"
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("D45").Value = "1.777.66"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "0.0₈110"
$ws.Range("E46").Value = "  -8.95%  "
$ws.Range("D47").Value = "'0.4439"
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("D49").Value = "'0.05056"
$ws.Range("E49").Value = "  -3.09%  "
$ws.Range("D50").Value = "'7.601"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("E51").Value = "  -0.52%  "
